$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 95, shifting existing rows (95..128) down to (96..129).
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new weekly data point.
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 44825
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = 100114002
$ws.Range("G95").Value = "Camote"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 50
$ws.Range("K95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("M95").Value = 20000
$ws.Range("N95").Value = "$/malla 20 kilos"
$ws.Range("O95").Value = "Perú"
$ws.Range("P95").Value = 1000
$ws.Range("Q95").Value = 20
$ws.Range("R95").Value = "Hortaliza"
